$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.930.42"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.22%  '
$ws.Range("D3").Value = "'2.582.91"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.58%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").Value = "'521.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.03%  '
$ws.Range("D6").Value = "'138.79"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.74%  '
$ws.Range("D7").Value = "'0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.10%  '
$ws.Range("D8").Value = "'0.564"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.63%  '
$ws.Range("D9").Value = "'2.592.72"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.98%  '
$ws.Range("D10").Value = "'6.51"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.45%  '
$ws.Range("E11").Value = '  -0.72%  '
$ws.Range("D12").Value = "'0.329"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.69%  '
$ws.Range("D13").Value = "'0.134"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.93%  '
$ws.Range("D14").Value = "'3.044.97"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.39%  '
$ws.Range("D15").Value = "'58.869.66"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.19%  '
$ws.Range("D16").Value = "'20.54"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.22%  '
$ws.Range("D17").Value = "'2.590.25"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.81%  '
$ws.Range("E18").Value = '  -1.50%  '
$ws.Range("D19").Value = "'338.99"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.97%  '
$ws.Range("D20").Value = "'4.30"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.24%  '
$ws.Range("D21").Value = "'10.05"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.88%  '
$ws.Range("D22").Value = "'6.43"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.79%  '
$ws.Range("D23").Value = "'0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.07%  '
$ws.Range("D24").Value = "'66.04"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.52%  '
$ws.Range("D25").Value = "'0.167"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.10%  '
$ws.Range("D26").Value = "'0.402"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.36%  '
$ws.Range("D27").Value = "'0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.20%  '
$ws.Range("D28").Value = "'7.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.66%  '
$ws.Range("D30").Value = "'0.0₃0719"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.29%  '
$ws.Range("D31").Value = "'5.88"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -6.04%  '
$ws.Range("E32").Value = '  -0.66%  '
$ws.Range("D33").Value = "'18.65"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.15%  '
$ws.Range("D34").Value = "'149.46"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.18%  '
$ws.Range("D35").Value = "'3.96"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.00%  '
$ws.Range("E36").Value = '  -3.03%  '
$ws.Range("D37").Value = "'36.62"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.39%  '
$ws.Range("E38").Value = '  -0.32%  '
$ws.Range("D39").Value = "'0.822"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.72%  '
$ws.Range("D40").Value = "'0.810"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -6.93%  '
$ws.Range("D41").Value = "'3.51"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.25%  '
$ws.Range("D42").Value = "'0.998"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.19%  '
$ws.Range("D43").Value = "'271.59"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.09%  '
$ws.Range("D44").Value = "'10.74"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.00%  '
$ws.Range("D45").Value = "'0.591"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.28%  '
$ws.Range("D46").Value = "'0.0949"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.71%  '
$ws.Range("D47").Value = "'0.0515"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.75%  '
$ws.Range("D48").Value = "'18.33"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.06%  '
$ws.Range("D49").Value = "'1.966.14"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.27%  '
$ws.Range("E50").Value = '  -1.13%  '
$ws.Range("D51").Value = "'4.49"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.81%  '
